$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 7, pushing existing rows 7-18 down to 8-19.
$ws.Rows.Item(7).Insert()

# Populate the newly inserted row 7 with the new weekly data point.
# All non-numeric/descriptive columns are identical to the other rows
# for this market/category combination.
$ws.Range("A7").Value = 1
$ws.Range("B7").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C7").Value = "Arica y Parinacota"
$ws.Range("D7").Value = 44679
$ws.Range("E7").Value = 15
$ws.Range("F7").Value = 100112003
$ws.Range("G7").Value = "Ajo"
$ws.Range("H7").Value = "Chino"
$ws.Range("I7").Value = "Primera"
$ws.Range("J7").Value = 200
$ws.Range("K7").Value = 19000
$ws.Range("L7").Value = 20000
$ws.Range("M7").Value = 19500
$ws.Range("N7").Value = "$/caja 10 kilos"
$ws.Range("O7").Value = "China"
$ws.Range("P7").Value = 1950
$ws.Range("Q7").Value = 10
$ws.Range("R7").Value = "Hortaliza"
